$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New text in the order that reproduces the shared-string table order
# ---------------------------------------------------------------------
$ws.Range("D4").Value = "Join date"
$ws.Range("B18").Value = "a"
$ws.Range("E4").Value = "IS GOOD"
$ws.Range("E6").Value = "'true"
$ws.Range("E5").Value = "'1"
$ws.Range("B7").Value = "Haha"
$ws.Range("E7").Value = "NO"

# ---------------------------------------------------------------------
# Row ids
# ---------------------------------------------------------------------
$ws.Range("A5").Value = 1
$ws.Range("A6").Value = 2
$ws.Range("A7").Value = 3

# ---------------------------------------------------------------------
# "Join date" column: D5 / D6 keep the stale quote-prefixed date style
# (as if the cell used to hold a typed '1' before the date was pasted in)
# ---------------------------------------------------------------------
$ws.Range("Z1").Value = "'1"
$ws.Range("Z1").NumberFormat = "mm-dd-yy"
$ws.Range("D5").Value = 45857
$ws.Range("Z1").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D6").Value = 42935
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("Z1").Clear()
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# C7 / D7: clean date style, reused from the existing C6 style (no
# quote-prefix baggage)
# ---------------------------------------------------------------------
$ws.Range("C6").Copy()
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("C7").Value = 42935
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = 42935
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Drop the now unused rows, relocate "aaaa" -> "a" into row 10
# ---------------------------------------------------------------------
$ws.Range("B9").Value = ""
$ws.Range("B15").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("B23").Value = ""
$ws.Range("B10").Value = "a"

# ---------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------
$ws.Columns("C").ColumnWidth = 10.7109375
$ws.Columns("D").ColumnWidth = 22.85546875

$ws.Range("D12").Select()
